$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 31 (shifts rows 31..96 down to 32..97) for the new
# "Email do fornecedor" field in the tb_fornecedor section.
$ws.Rows.Item(31).Insert()

# Copy formatting from the row above (row 30, an existing "tb_fornecedor"
# field row) onto the freshly inserted row 31 so borders/style match -
# restrict to the used A:D columns to avoid bloating the sheet dimension.
$ws.Range("A30:D30").Copy()
$ws.Range("A31:D31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the new row's values.
$ws.Range("A31").Value = "email"
$ws.Range("B31").Value = "varchar(50)"
$ws.Range("C31").Value = "Email do fornecedor"
$ws.Range("D31").Value = "not null; unique"

# Update the sheet view state (scroll position + active selection) to match
# the saved workbook state.
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("C31").Select()
